# Update 16 Jan 2023
# Add three new Area/Region mapping rows to Sheet1 (rows 113-115).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$newRows = @(
    @{ Row = 113; A = "West Sulawesi Province";      B = "Sulawesi" },
    @{ Row = 114; A = "Riau islands";                B = "Sumatera 1" },
    @{ Row = 115; A = "West Nusa Tenggara Province";  B = "BARA" }
)

# Copy the formatting of the last existing data row (A112) onto the new
# column-A cells so the new rows keep the same look (font/alignment) as
# the rest of the table, without introducing new style/font entries.
$ws.Range("A112").Copy() | Out-Null
foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = $false

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
}

# Restore the view state (selection sits just past the new last row).
$ws.Range("B119").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
